$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.Value = "'302.48"
$r.Style = "Normal"

$r = $ws.Range("E2")
$r.Value = "'0.74%"
$r.Style = "Normal"

$r = $ws.Range("G2")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D3")
$r.Value = "'32.16"
$r.Style = "Normal"

$r = $ws.Range("E3")
$r.Value = "'1.49%"
$r.Style = "Normal"

$r = $ws.Range("G3")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D4")
$r.Value = "'4.987"
$r.Style = "Normal"

$r = $ws.Range("E4")
$r.Value = "'-3.15%"
$r.Style = "Normal"

$r = $ws.Range("G4")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D5")
$r.Value = "'0.07907"
$r.Style = "Normal"

$r = $ws.Range("E5")
$r.Value = "'-2.62%"
$r.Style = "Normal"

$r = $ws.Range("G5")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D6")
$r.Value = "'2.102"
$r.Style = "Normal"

$r = $ws.Range("E6")
$r.Value = "'-16.47%"
$r.Style = "Normal"

$r = $ws.Range("G6")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D7")
$r.Value = "'7.834"
$r.Style = "Normal"

$r = $ws.Range("E7")
$r.Value = "'0.59%"
$r.Style = "Normal"

$r = $ws.Range("G7")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D8")
$r.Value = "'3.802"
$r.Style = "Normal"

$r = $ws.Range("E8")
$r.Value = "'-2.22%"
$r.Style = "Normal"

$r = $ws.Range("G8")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D9")
$r.Value = "'0.9263"
$r.Style = "Normal"

$r = $ws.Range("E9")
$r.Value = "'0.10%"
$r.Style = "Normal"

$r = $ws.Range("G9")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D10")
$r.Value = "'0.1743"
$r.Style = "Normal"

$r = $ws.Range("E10")
$r.Value = "'-1.06%"
$r.Style = "Normal"

$r = $ws.Range("G10")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D11")
$r.Value = "'0.07942"
$r.Style = "Normal"

$r = $ws.Range("E11")
$r.Value = "'6.32%"
$r.Style = "Normal"

$r = $ws.Range("G11")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D12")
$r.Value = "'0.08608"
$r.Style = "Normal"

$r = $ws.Range("E12")
$r.Value = "'-2.36%"
$r.Style = "Normal"

$r = $ws.Range("G12")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D13")
$r.Value = "'0.03134"
$r.Style = "Normal"

$r = $ws.Range("E13")
$r.Value = "'3.55%"
$r.Style = "Normal"

$r = $ws.Range("G13")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D14")
$r.Value = "'0.1003"
$r.Style = "Normal"

$r = $ws.Range("E14")
$r.Value = "'0.29%"
$r.Style = "Normal"

$r = $ws.Range("G14")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D15")
$r.Value = "'0.001522"
$r.Style = "Normal"

$r = $ws.Range("E15")
$r.Value = "'0.87%"
$r.Style = "Normal"

$r = $ws.Range("G15")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D16")
$r.Value = "'0.005802"
$r.Style = "Normal"

$r = $ws.Range("E16")
$r.Value = "'-0.30%"
$r.Style = "Normal"

$r = $ws.Range("G16")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("E17")
$r.Value = "'2,098.28%"
$r.Style = "Normal"

$r = $ws.Range("G17")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D18")
$r.Value = "'3.466"
$r.Style = "Normal"

$r = $ws.Range("E18")
$r.Value = "'-2.93%"
$r.Style = "Normal"

$r = $ws.Range("G18")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("E19")
$r.Value = "'-0.43%"
$r.Style = "Normal"

$r = $ws.Range("G19")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("E20")
$r.Value = "'-0.07%"
$r.Style = "Normal"

$r = $ws.Range("G20")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D21")
$r.Value = "'0.1291"
$r.Style = "Normal"

$r = $ws.Range("E21")
$r.Value = "'-3.67%"
$r.Style = "Normal"

$r = $ws.Range("G21")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D22")
$r.Value = "'4.322"
$r.Style = "Normal"

$r = $ws.Range("E22")
$r.Value = "'3.93%"
$r.Style = "Normal"

$r = $ws.Range("G22")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D23")
$r.Value = "'0.1790"
$r.Style = "Normal"

$r = $ws.Range("E23")
$r.Value = "'6.45%"
$r.Style = "Normal"

$r = $ws.Range("G23")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D24")
$r.Value = "'0.04598"
$r.Style = "Normal"

$r = $ws.Range("E24")
$r.Value = "'-0.48%"
$r.Style = "Normal"

$r = $ws.Range("G24")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D25")
$r.Value = "'0.001237"
$r.Style = "Normal"

$r = $ws.Range("E25")
$r.Value = "'-0.31%"
$r.Style = "Normal"

$r = $ws.Range("G25")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D26")
$r.Value = "'0.004473"
$r.Style = "Normal"

$r = $ws.Range("G26")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D27")
$r.Value = "'0.0001250"
$r.Style = "Normal"

$r = $ws.Range("E27")
$r.Value = "'4.11%"
$r.Style = "Normal"

$r = $ws.Range("G27")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("G28")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("G29")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("G30")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("G31")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("G32")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("G33")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("G34")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("G35")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("G36")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("G37")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("G38")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D39")
$r.Value = "'0.01722"
$r.Style = "Normal"

$r = $ws.Range("G39")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D40")
$r.Value = "'0.04792"
$r.Style = "Normal"

$r = $ws.Range("E40")
$r.Value = "'4.07%"
$r.Style = "Normal"

$r = $ws.Range("G40")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D41")
$r.Value = "'0.007467"
$r.Style = "Normal"

$r = $ws.Range("E41")
$r.Value = "'7.26%"
$r.Style = "Normal"

$r = $ws.Range("G41")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D42")
$r.Value = "'0.1360"
$r.Style = "Normal"

$r = $ws.Range("E42")
$r.Value = "'-0.99%"
$r.Style = "Normal"

$r = $ws.Range("G42")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D43")
$r.Value = "'0.002390"
$r.Style = "Normal"

$r = $ws.Range("E43")
$r.Value = "'9.07%"
$r.Style = "Normal"

$r = $ws.Range("G43")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D44")
$r.Value = "'0.01023"
$r.Style = "Normal"

$r = $ws.Range("E44")
$r.Value = "'-2.08%"
$r.Style = "Normal"

$r = $ws.Range("G44")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D45")
$r.Value = "'0.00005975"
$r.Style = "Normal"

$r = $ws.Range("E45")
$r.Value = "'-3.64%"
$r.Style = "Normal"

$r = $ws.Range("G45")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D46")
$r.Value = "'0.00000000750"
$r.Style = "Normal"

$r = $ws.Range("E46")
$r.Value = "'-0.04%"
$r.Style = "Normal"

$r = $ws.Range("G46")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D47")
$r.Value = "'0.003391"
$r.Style = "Normal"

$r = $ws.Range("E47")
$r.Value = "'-59.66%"
$r.Style = "Normal"

$r = $ws.Range("G47")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("E48")
$r.Value = "'2.75%"
$r.Style = "Normal"

$r = $ws.Range("G48")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D49")
$r.Value = "'0.00002101"
$r.Style = "Normal"

$r = $ws.Range("E49")
$r.Value = "'-0.04%"
$r.Style = "Normal"

$r = $ws.Range("G49")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("D50")
$r.Value = "'0.0002001"
$r.Style = "Normal"

$r = $ws.Range("E50")
$r.Value = "'-0.04%"
$r.Style = "Normal"

$r = $ws.Range("G50")
$r.Value = "'22"
$r.Style = "Normal"

$r = $ws.Range("G51")
$r.Value = "'22"
$r.Style = "Normal"
